# Adds two new worksheets ("1_Thinking" and "Sheet2") after "0_Bears",
# populates "1_Thinking" with an inductive-reasoning question (with
# word-wrapped cells, custom column widths and row heights), makes
# "Sheet2" the active/selected sheet, and updates the leftover
# selection on "0_Bears".

$wb  = $excel.ActiveWorkbook
$ws0 = $wb.Worksheets.Item(1)

# --- create the new sheets, in order, right after "0_Bears" ---------------
$wsThinking = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws0)
$wsThinking.Name = "1_Thinking"

$wsSheet2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsThinking)
$wsSheet2.Name = "Sheet2"

# --- populate "1_Thinking" (write column-by-column so shared-string
#     interning order matches the authored workbook) -----------------------
$wsThinking.Range("A1").Value = "`"Inductive Reasoning`" is often described as going `"from specifics to the general`": observing individual facts and trying to come up with a general rule about how something works.   Is this closer to `"thinking like a human`" or `"thinking like a computer`"?"
$wsThinking.Range("A2").Value = "Human"
$wsThinking.Range("A3").Value = "Computer"

$wsThinking.Range("B1").Value = "Correct"
$wsThinking.Range("B2").Value = "Y"
$wsThinking.Range("B3").Value = "N"

$wsThinking.Range("C1").Value = "Comment"
$wsThinking.Range("C2").Value = "Inductive reasoning is a kind of modeling: the `"general rule`" that it seeks is essentially a model."
$wsThinking.Range("C3").Value = "Computers usually start with a rule (an algorithm) and decide whether the data fits that rule (`"if *this* is greater than *that*, then do this….`")"

# column widths
$wsThinking.Columns.Item(1).ColumnWidth = 38.42578125
$wsThinking.Columns.Item(2).ColumnWidth = 31.85546875
$wsThinking.Columns.Item(3).ColumnWidth = 31.42578125

# wrap text on the data range, with row heights matching the wrapped content
$wsThinking.Range("A1:C3").WrapText = $true
$wsThinking.Rows.Item(1).RowHeight = 105
$wsThinking.Rows.Item(2).RowHeight = 45
$wsThinking.Rows.Item(3).RowHeight = 75

$wsThinking.Range("C4").Select() | Out-Null

# --- restore/update the selection left on "0_Bears" -------------------------
$ws0.Range("A1:C7").Select() | Out-Null

# --- "Sheet2" is the sheet left active/selected in the saved workbook ------
$wsSheet2.Activate() | Out-Null
